# Weekly update: insert two new daily-price records for Ciboulette
# (Mercado Mayorista Lo Valledor de Santiago) at the top of the
# historical block (row 423), pushing the existing data down by two
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 423 (inserting twice at the same index
# pushes the prior insert down, giving us two fresh rows at 423/424
# and shifting everything that used to start at row 423 down to 425).
$ws.Rows.Item(423).Insert()
$ws.Rows.Item(423).Insert()

# New row 423: Primera quality
$ws.Cells.Item(423,1).Value = 6
$ws.Cells.Item(423,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(423,3).Value = "Metropolitana"
$ws.Cells.Item(423,4).Value = 44798
$ws.Cells.Item(423,5).Value = 13
$ws.Cells.Item(423,6).Value = 100112039
$ws.Cells.Item(423,7).Value = "Ciboulette"
$ws.Cells.Item(423,8).Value = "Sin especificar"
$ws.Cells.Item(423,9).Value = "Primera"
$ws.Cells.Item(423,10).Value = 580
$ws.Cells.Item(423,11).Value = 1800
$ws.Cells.Item(423,12).Value = 2000
$ws.Cells.Item(423,13).Value = 1879
$ws.Cells.Item(423,14).Value = "`$/docena de atados"
$ws.Cells.Item(423,15).Value = "Región Metropolitana"
$ws.Cells.Item(423,16).Value = 626
$ws.Cells.Item(423,17).Value = 3
$ws.Cells.Item(423,18).Value = "Hortaliza"

# New row 424: Segunda quality
$ws.Cells.Item(424,1).Value = 6
$ws.Cells.Item(424,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(424,3).Value = "Metropolitana"
$ws.Cells.Item(424,4).Value = 44798
$ws.Cells.Item(424,5).Value = 13
$ws.Cells.Item(424,6).Value = 100112039
$ws.Cells.Item(424,7).Value = "Ciboulette"
$ws.Cells.Item(424,8).Value = "Sin especificar"
$ws.Cells.Item(424,9).Value = "Segunda"
$ws.Cells.Item(424,10).Value = 200
$ws.Cells.Item(424,11).Value = 1300
$ws.Cells.Item(424,12).Value = 1300
$ws.Cells.Item(424,13).Value = 1300
$ws.Cells.Item(424,14).Value = "`$/docena de atados"
$ws.Cells.Item(424,15).Value = "Región Metropolitana"
$ws.Cells.Item(424,16).Value = 433
$ws.Cells.Item(424,17).Value = 3
$ws.Cells.Item(424,18).Value = "Hortaliza"
